# Apply updated registration counts (Inscritos / Pagos / Inscrições homologadas)
# to the "Inscricoes" worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row -> column -> new value
$updates = @{
    2  = @{ E = 82;  F = 57;  H = 57  }
    10 = @{ E = 353 }
    11 = @{ E = 241 }
    12 = @{ E = 350; F = 193; H = 193 }
    13 = @{ E = 102 }
    17 = @{ E = 66  }
    21 = @{ E = 115; F = 59;  H = 59  }
    23 = @{ E = 158 }
    24 = @{ E = 156; F = 76;  H = 76  }
    25 = @{ E = 183; F = 84;  H = 84  }
    29 = @{ E = 139 }
    30 = @{ E = 159; F = 88;  H = 88  }
    32 = @{ E = 142; F = 76;  H = 76  }
    33 = @{ E = 230 }
    38 = @{ E = 76  }
    40 = @{ E = 204 }
    41 = @{ E = 297; F = 127; H = 127 }
    42 = @{ E = 263 }
    43 = @{ E = 87  }
    45 = @{ E = 101; F = 44;  H = 44  }
    46 = @{ E = 229 }
    47 = @{ E = 332 }
    48 = @{ E = 151 }
    49 = @{ E = 219 }
    50 = @{ E = 189; F = 67;  H = 67  }
}

foreach ($row in $updates.Keys) {
    $cols = $updates[$row]
    foreach ($col in $cols.Keys) {
        $cellAddr = "$col$row"
        $ws.Range($cellAddr).Value = $cols[$col]
    }
}
